$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark after "Check list for COA.COV"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Fix typo "R-16)" -> "R-61)" and re-insert the "_GoBack" bookmark
#    right after the new "R-61" text (before the closing paren).
$rng = $d.Content
$rng.Find.Execute("R-16", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "R-61"
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)

# 3) Merge "Complete info and" / " check "Motion"" into a single run
$rng = $d.Content
$rng.Find.Execute("and check", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "andXXXcheck"
$rng = $d.Content
$rng.Find.Execute("andXXXcheck", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "and check"

# 4) Merge "Add names of client under "In the Matters o" / "f" " into one run
$rng = $d.Content
$rng.Find.Execute("of”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "ofXXX”"
$rng = $d.Content
$rng.Find.Execute("ofXXX”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "of”"

# 5) Merge "...sign on their beh" / "alf. " into a single run
$rng = $d.Content
$rng.Find.Execute("behalf. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "behalfXXX. "
$rng = $d.Content
$rng.Find.Execute("behalfXXX. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "behalf. "
